$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to carry a hidden helper column (B) driving an AutoFilter
# on B3 (and its associated _FilterDatabase defined name). Turn the filter
# off and drop the defined name it leaves behind, then delete the now
# unneeded helper column so C:F shifts left into B:E and is no longer
# hidden.
$ws.AutoFilterMode = $false
foreach ($n in $wb.Names) {
    $n.Delete()
}
$ws.Columns.Item(2).Delete()

# Fill in the Origen/Destino sample rows requested by the update.
$ws.Range("B4").Value = "Bogota "
$ws.Range("C4").Value = "Cali "

# Put the active selection on column B like the refreshed template.
$ws.Range("B1:B1048576").Select()
